# Auto-generated script applying cryptos.xlsx price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices/percentages) are kept as exact text
# by forcing the Text number format on the cells we are about to write, matching
# the original inline-string/text representation used in the workbook.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "74.847.15"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "2.812.57"
$ws.Range("E3").Value = "  +7.38%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "187.73"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "591.57"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("D9").Value = "0.188"
$ws.Range("E9").Value = "  -3.98%  "
$ws.Range("D10").Value = "2.811.69"
$ws.Range("E10").Value = "  +7.51%  "
$ws.Range("D11").Value = "0.375"
$ws.Range("E11").Value = "  +5.38%  "
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").Value = "3.326.80"
$ws.Range("E14").Value = "  +7.21%  "
$ws.Range("D15").Value = "74.732.84"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "0.0000186"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "26.82"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "2.810.77"
$ws.Range("E18").Value = "  +7.64%  "
$ws.Range("D19").Value = "8.87"
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").Value = "12.25"
$ws.Range("E20").Value = "  +4.19%  "
$ws.Range("D21").Value = "376.24"
$ws.Range("E21").Value = "  +3.11%  "
$ws.Range("D22").Value = "2.26"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "70.61"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "2.948.34"
$ws.Range("E26").Value = "  +7.08%  "
$ws.Range("D27").Value = "4.14"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "9.66"
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0000102"
$ws.Range("E29").Value = "  +9.18%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "1.39"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "509.72"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("D33").Value = "7.62"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "164.58"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").Value = "19.82"
$ws.Range("E37").Value = "  +4.00%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "19.35"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "182.02"
$ws.Range("E40").Value = "  +12.68%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "0.340"
$ws.Range("E42").Value = "  +4.78%  "
$ws.Range("D43").Value = "4.98"
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("D44").Value = "1.66"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "1.20"
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "39.97"
$ws.Range("E46").Value = "  +2.59%  "
$ws.Range("D47").Value = "0.0863"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").Value = "2.30"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").Value = "0.574"
$ws.Range("E49").Value = "  +10.13%  "
$ws.Range("D50").Value = "3.71"
$ws.Range("E50").Value = "  +3.28%  "
$ws.Range("D51").Value = "0.631"
$ws.Range("E51").Value = "  +8.08%  "
